$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# New rows of OBI problems, entering cell values in the exact order that
# reproduces the shared-string insertion order seen in the target file.
$ws.Range("B4").Value = "Soma"
$ws.Range("D5").Value = "função recursiva"
$ws.Range("B5").Value = "Fatorial"
$ws.Range("C4").Value = "SOMA"
$ws.Range("C5").Value = "FATORIA2"
$ws.Range("B6").Value = "Quermesse"
$ws.Range("C6").Value = "QUERM"

$ws.Range("A4").Value = 3830
$ws.Range("D4").Value = "-"

$ws.Range("A5").Value = 3774

$ws.Range("A6").Value = 811
$ws.Range("D6").Value = "-"

$ws.Range("B11").Select() | Out-Null

$wb.Save()
